$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily mod-count entry as row 36.
# The date is entered with a leading apostrophe so it is stored as literal
# text ("2025/12/15") instead of being auto-converted to a date serial,
# matching the existing column A values which are plain text dates.
$ws.Range("A36").Value = "'2025/12/15"
$ws.Range("B36").Value = "逃离鸭科夫"
$ws.Range("C36").Value = 1348

# Match the centered alignment used by the rest of the data rows (row 2-35).
$ws.Range("A36:C36").HorizontalAlignment = -4108
$ws.Range("A36:C36").VerticalAlignment = -4108
